$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation (leading apostrophe) to preserve values such as
# "8.0", "10", and thousand-separated numbers as literal text, matching the
# original inlineStr cell contents instead of being auto-converted to numbers.

$ws.Range('B2').Value = '''US$1,116'
$ws.Range('A3').Value = '''Sonder Le Frochot'
$ws.Range('B3').Value = '''US$682'
$ws.Range('C3').Value = '''7.5'
$ws.Range('D3').Value = '''Good'
$ws.Range('E3').Value = '''724'
$ws.Range('A5').Value = '''St Christopher''s Inn Paris - Canal'
$ws.Range('B5').Value = '''US$86'
$ws.Range('C5').Value = '''7.3'
$ws.Range('E5').Value = '''6,519'
$ws.Range('A6').Value = '''The People - Paris Nation'
$ws.Range('B6').Value = '''US$121'
$ws.Range('C6').Value = '''8.4'
$ws.Range('D6').Value = '''Very Good'
$ws.Range('E6').Value = '''6,122'
$ws.Range('A7').Value = '''Hotel Relais Bosquet by Malone'
$ws.Range('B7').Value = '''US$924'
$ws.Range('E7').Value = '''1,068'
$ws.Range('A8').Value = '''Hotel Armoni Paris'
$ws.Range('B8').Value = '''US$550'
$ws.Range('C8').Value = '''8.1'
$ws.Range('D8').Value = '''Very Good'
$ws.Range('E8').Value = '''1,899'
$ws.Range('A9').Value = '''Odalys City Paris XVII'
$ws.Range('B9').Value = '''US$571'
$ws.Range('C9').Value = '''7.5'
$ws.Range('E9').Value = '''3,845'
$ws.Range('A10').Value = '''St Christopher''s Inn Paris - Gare du Nord'
$ws.Range('B10').Value = '''US$95'
$ws.Range('C10').Value = '''7.2'
$ws.Range('D10').Value = '''Good'
$ws.Range('E10').Value = '''9,248'
$ws.Range('A11').Value = '''Auberge de Jeunesse HI Paris Yves Robert'
$ws.Range('B11').Value = '''US$127'
$ws.Range('C11').Value = '''7.7'
$ws.Range('D11').Value = '''Good'
$ws.Range('E11').Value = '''5,517'
$ws.Range('A12').Value = '''Hotel The Playce by Happyculture'
$ws.Range('B12').Value = '''US$450'
$ws.Range('C12').Value = '''7.7'
$ws.Range('D12').Value = '''Good'
$ws.Range('E12').Value = '''2,242'
$ws.Range('A13').Value = '''Le Regent Montmartre by Hiphophostels'
$ws.Range('B13').Value = '''US$99'
$ws.Range('C13').Value = '''7.4'
$ws.Range('D13').Value = '''Good'
$ws.Range('E13').Value = '''5,885'
$ws.Range('A14').Value = '''citizenM Paris Champs-Élysées'
$ws.Range('B14').Value = '''US$1,167'
$ws.Range('C14').Value = '''8.5'
$ws.Range('D14').Value = '''Very Good'
$ws.Range('E14').Value = '''2,787'
$ws.Range('A15').Value = '''Motel One Paris-Porte Dorée'
$ws.Range('B15').Value = '''US$495'
$ws.Range('C15').Value = '''8.7'
$ws.Range('D15').Value = '''Excellent'
$ws.Range('E15').Value = '''6,979'
$ws.Range('A16').Value = '''The People - Paris Bercy'
$ws.Range('B16').Value = '''US$143'
$ws.Range('C16').Value = '''8.8'
$ws.Range('D16').Value = '''Excellent'
$ws.Range('E16').Value = '''5,976'
$ws.Range('A17').Value = '''B&B HOTEL Paris 17 Batignolles'
$ws.Range('B17').Value = '''US$506'
$ws.Range('C17').Value = '''7.9'
$ws.Range('E17').Value = '''12,198'
$ws.Range('A18').Value = '''Beau M Paris'
$ws.Range('B18').Value = '''US$138'
$ws.Range('C18').Value = '''8.7'
$ws.Range('D18').Value = '''Excellent'
$ws.Range('E18').Value = '''1,149'
$ws.Range('A19').Value = '''Hôtel Cabane - Orso Hotels'
$ws.Range('B19').Value = '''US$528'
$ws.Range('C19').Value = '''8.4'
$ws.Range('D19').Value = '''Very Good'
$ws.Range('E19').Value = '''983'
$ws.Range('A20').Value = '''Hôtel le 209 Paris Bercy'
$ws.Range('B20').Value = '''US$547'
$ws.Range('C20').Value = '''8.3'
$ws.Range('D20').Value = '''Very Good'
$ws.Range('E20').Value = '''3,245'
$ws.Range('A21').Value = '''Hôtel Le Daum'
$ws.Range('B21').Value = '''US$536'
$ws.Range('C21').Value = '''8.4'
$ws.Range('E21').Value = '''943'
$ws.Range('A22').Value = '''ibis Paris Gare du Nord Château Landon 10ème'
$ws.Range('B22').Value = '''US$551'
$ws.Range('C22').Value = '''7.5'
$ws.Range('D22').Value = '''Good'
$ws.Range('E22').Value = '''4,930'
$ws.Range('A23').Value = '''Generator Paris'
$ws.Range('B23').Value = '''US$103'
$ws.Range('C23').Value = '''7.3'
$ws.Range('E23').Value = '''9,523'
$ws.Range('A24').Value = '''Elysees Niel Hotel'
$ws.Range('B24').Value = '''US$510'
$ws.Range('C24').Value = '''8.0'
$ws.Range('D24').Value = '''Very Good'
$ws.Range('E24').Value = '''1,369'
$ws.Range('A25').Value = '''Hôtel Belloy Saint Germain'
$ws.Range('B25').Value = '''US$776'
$ws.Range('C25').Value = '''8.0'
$ws.Range('D25').Value = '''Very Good'
$ws.Range('E25').Value = '''718'
$ws.Range('A26').Value = '''Hotel Du Cadran'
$ws.Range('B26').Value = '''US$784'
$ws.Range('C26').Value = '''8.5'
$ws.Range('E26').Value = '''795'
$ws.Range('A27').Value = '''Hotel Monceau Wagram'
$ws.Range('B27').Value = '''US$599'
$ws.Range('C27').Value = '''7.3'
$ws.Range('D27').Value = '''Good'
$ws.Range('E27').Value = '''1,379'
